# Ex18 Initio Simulator and tidying part 3.
#
# 1. Title line: "Pi2Go Simulator Programming: " -> "Virtual Pi2Go Programming: "
# 2. Exercise heading: "Exercise:  " -> "Exercise 1:  " (with _GoBack bookmark
#    relocated to just before the trailing ":  ")

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: Title line
# ---------------------------------------------------------------------------

# Remove "Simulator " - it merges with the following "Programming: " run,
# producing a single run that reads " Programming: ".
$full = $d.Content.Text
$idxSim = $full.IndexOf("Simulator ")
$rSim = $d.Range($idxSim, $idxSim + 10)
$rSim.Delete()

# Turn the leading "Pi2Go" run into "Virtual Pi2Go" (still a single run).
$rTitle = $d.Range(0, 0)
$rTitle.Text = "Virtual Pi2Go"

# Force a clean run split between "Virtual " and "Pi2Go" by dropping a
# temporary bookmark at the boundary and then removing it again - the
# bookmark forces the run boundary without leaving any formatting residue.
$rBoundary = $d.Range(8, 8)
$d.Bookmarks.Add("zzzTmpTitleSplit", $rBoundary)
$d.Bookmarks("zzzTmpTitleSplit").Delete()

# ---------------------------------------------------------------------------
# Part 2: "Exercise:  " -> "Exercise" / " 1" / ":  " with the _GoBack bookmark
# relocated in between the new " 1" run and the trailing ":  " run.
# ---------------------------------------------------------------------------

$full = $d.Content.Text
$idxEx = $full.IndexOf("Exercise:")
$insPt = $idxEx + 8   # just after "Exercise"
$rIns = $d.Range($insPt, $insPt)
$rIns.InsertBefore(" 1")

$splitA = $idxEx + 8       # boundary between "Exercise" and " 1"
$splitB = $idxEx + 8 + 2   # boundary between " 1" and ":  "

# Force the "Exercise" | " 1" boundary with a temporary bookmark.
$rSplitA = $d.Range($splitA, $splitA)
$d.Bookmarks.Add("zzzTmpExSplit", $rSplitA)

# Relocate (or create) the real "_GoBack" bookmark at the " 1" | ":  "
# boundary - Bookmarks.Add moves an existing bookmark of the same name,
# which also removes it from its old location next to
# "When testing your program,".
$rSplitB = $d.Range($splitB, $splitB)
$d.Bookmarks.Add("_GoBack", $rSplitB)

$d.Bookmarks("zzzTmpExSplit").Delete()
